$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-Contact($row, $first, $last, $email, $phone) {
    $ws.Cells.Item($row, 1).Value = $first
    $ws.Cells.Item($row, 2).Value = $last
    $ws.Cells.Item($row, 3).Value = $email
    $ws.Cells.Item($row, 4).Value = $phone
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 3), "mailto:" + $email) | Out-Null
    $ws.Cells.Item($row, 3).Style = $ws.Range("C1").Style
}

# --- Append new contact rows (rows 19-35) ---
Add-Contact 19 "Ethan" "Carter" "ecarter@yahoo.com" "0957642231"
Add-Contact 20 "Olivia" "Hartman" "ohartman@live.com" "0981616544"
Add-Contact 21 "Julian" "Mercer" "jmercer@live.com" "0998172743"
Add-Contact 22 "Chloe" "Davenport" "cdavenport@gmail.com" "0917583245"
Add-Contact 23 "Gavin" "Halstead" "ghalstead@yahoo.com" "0927653262"
Add-Contact 24 "Lila" "Kensington" "lkensington@gmail.com" "0921345354"
Add-Contact 25 "Daniel" "Rowley" "drowley@live.com" "0988787977"
Add-Contact 26 "Harper" "Linton" "hlinton@gmail.com" "0972134890"
Add-Contact 27 "Xavier" "Callahan" "xcallahan@gmail.com" "0918765321"

# Rows 28-29 (Zoe Merritt / Tristan Holloway) were originally typed
# out of column order: both names + row 28's phone first, then both
# emails, then row 29's phone - reproduce that entry order here.
$ws.Cells.Item(28, 1).Value = "Zoe"
$ws.Cells.Item(28, 2).Value = "Merritt"
$ws.Cells.Item(28, 4).Value = "0996876644"
$ws.Cells.Item(29, 1).Value = "Tristan"
$ws.Cells.Item(29, 2).Value = "Holloway"
$ws.Cells.Item(28, 3).Value = "zmerritt@yahoo.com"
$ws.Cells.Item(29, 3).Value = "tholloway@live.com"
$ws.Cells.Item(29, 4).Value = "0957987981"

$ws.Hyperlinks.Add($ws.Cells.Item(29, 3), "mailto:tholloway@live.com") | Out-Null
$ws.Cells.Item(29, 3).Style = $ws.Range("C1").Style
$ws.Hyperlinks.Add($ws.Cells.Item(28, 3), "mailto:zmerritt@yahoo.com") | Out-Null
$ws.Cells.Item(28, 3).Style = $ws.Range("C1").Style

Add-Contact 30 "Maya" "Pennington" "mpennington@yahoo.com" "0918374653"
Add-Contact 31 "Nora" "Becket" "nbecket@live.com" "0986567214"
Add-Contact 32 "Spencer" "Aldridge" "saldridge@gmail.com" "0916783243"
Add-Contact 33 "Hazel" "Kingsley" "hkingsley@gmail.com" "0992345612"
Add-Contact 34 "Owen" "Thorne" "othorne@yahoo.com" "0916366327"
Add-Contact 35 "Amelia" "Prescott" "aprescott@gmail.com" "0928764378"

# --- Restore the view state (scroll position + selection) ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C38").Select()
